# "replaced with new main"
#
# Target text: the first paragraph currently reads
#     This is a Microsoft word document.
# and must become (as four *separate* runs, matching the author's diff):
#     [This is a Microsoft word document.][ (][Changed main][)]
#
# A plain Range.InsertAfter / Find.Execute replace ends up coalescing the
# inserted text into the existing run (same rPr => same run), which would
# leave a single merged run instead of the three extra <w:r> elements the
# diff shows. To preserve the run boundaries we splice in real OOXML via
# Range.InsertXML for the new runs (that keeps each <w:r> distinct), then
# remove the paragraph break it necessarily introduces so everything lands
# back in the original paragraph (carrying over the original paragraph's
# identity attributes, e.g. w14:paraId, so that the <w:p> tag itself is
# left untouched, exactly like the diff).

$d = $word.ActiveDocument

$firstPara = $d.Paragraphs(1).Range

# Capture the existing paragraph's opening <w:p ...> tag (with its
# w14:paraId / w:rsidR / etc.) so the re-merged paragraph keeps it intact.
$null = $firstPara.WordOpenXML -match '<w:p( [^>]*)?>'
$paraOpenTag = $matches[0]

# Position right before the paragraph mark, i.e. the true end of
# "This is a Microsoft word document."
$insertAt = $firstPara.End - 1
$target = $d.Range($insertAt, $insertAt)

$newParaXml = $paraOpenTag +
  '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
  '<w:r><w:t>Changed main</w:t></w:r>' +
  '<w:r><w:t>)</w:t></w:r>' +
  '</w:p>'

$package = '<?xml version="1.0" standalone="yes"?>' +
  '<?mso-application progid="Word.Document"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">' +
      '<pkg:xmlData>' +
        '<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' +
          '<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>' +
        '</Relationships>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
          '<w:body>' + $newParaXml + '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

# InsertXML always creates its content as a fresh paragraph, splitting the
# original paragraph in two. Immediately merge the split back together by
# deleting the paragraph mark between them; the surviving mark (and thus
# the <w:p> attributes) come from the new paragraph, which is why we
# seeded it above with the original paragraph's own opening tag.
$target.InsertXML($package)

$seam = $d.Range($firstPara.End - 1, $firstPara.End)
$seam.Delete()
